# Added handling of common packages.
#
# The "methodNumberOfLines" sheet gains a new row describing the number of
# lines in SecuritySecureConfig's constructor. It is inserted above the
# existing "configure(...)" row (which holds the same Class Name), pushing
# the later rows down by one.
#
# Row 3 already has the exact text we need in column C ("3" via the shared
# string used for the number-of-lines value), so duplicate that row upward
# (Copy + Insert, which preserves cell types/styles natively, unlike
# re-assigning .Value which would coerce a numeric-looking string like "3"
# into a literal number) and then only overwrite the two cells that must
# change (A and B); column C is already correct and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

$ws.Rows.Item(3).Copy() | Out-Null
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "com.macro.mall.config.SecuritySecureConfig"
$ws.Range("B2").Value = "SecuritySecureConfig(de.codecentric.boot.admin.server.config.AdminServerProperties)"
